$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '59.416.76'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.521.20'
$ws.Range('E3').Value = '  +1.05%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '543.38'
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.47'
$ws.Range('E6').Value = '  -1.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  -0.81%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.542.36'
$ws.Range('E9').Value = '  +0.89%  '
$ws.Range('E10').Value = '  -0.22%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.58'
$ws.Range('E12').Value = '  +1.55%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.361'
$ws.Range('E13').Value = '  +0.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.963.50'
$ws.Range('E14').Value = '  +0.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.54'
$ws.Range('E15').Value = '  -4.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '59.502.44'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('E17').Value = '  +1.06%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.529.68'
$ws.Range('E18').Value = '  +1.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.20'
$ws.Range('E19').Value = '  -2.19%  '
$ws.Range('E20').Value = '  -2.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '326.52'
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('E22').Value = '  +0.78%  '
$ws.Range('E23').Value = '  +1.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.60'
$ws.Range('E24').Value = '  +1.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.434'
$ws.Range('E25').Value = '  -3.16%  '
$ws.Range('E26').Value = '  +1.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.995'
$ws.Range('E27').Value = '  -1.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.01'
$ws.Range('E28').Value = '  +2.24%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0789'
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.73'
$ws.Range('E31').Value = '  -0.85%  '
$ws.Range('E32').Value = '  -8.62%  '
$ws.Range('E33').Value = '  +2.67%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '160.40'
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('B35').Value = 'USDe'
$ws.Range('C35').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  +0.13%  '
$ws.Range('E36').Value = '  -1.31%  '
$ws.Range('E37').Value = '  -2.27%  '
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('E40').Value = '  -8.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.836'
$ws.Range('E41').Value = '  +0.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '296.33'
$ws.Range('E42').Value = '  -5.82%  '
$ws.Range('E43').Value = '  -2.58%  '
$ws.Range('E44').Value = '  +0.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.603'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.83'
$ws.Range('E46').Value = '  +1.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0937'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.91'
$ws.Range('E48').Value = '  +1.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '123.32'
$ws.Range('E49').Value = '  -2.78%  '
$ws.Range('E50').Value = '  -1.99%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0514'
$ws.Range('E51').Value = '  -3.70%  '
